$d = $word.ActiveDocument

# ---- Step 1: rewrite paragraph 3 (the "3)" SCRUM paragraph) ----
$p3 = $d.Paragraphs.Item(3).Range
$pStart = $p3.Start

$run1 = "3)"
$run2 = " "
$run3old = "O SCRUM se adaptaria ao projeto, é um processo ágil com transparência, todos tem conhecimento dos requisitos e os processos e do andamento do processo, "
$run4old = "processo com constante inspeção, seja no sprint review ou nas reuniões diárias além de ser um processo extremamente adaptativo, utilizando este processo o cliente vai poder participar de perto do desenvolvimento, além de ter um rápido resultado, vai ajudar o cliente a decidir em qual ponto focar durante o desenvolvimento, já que ele tem a necessidade de muitos recursos."

$run3Start = $pStart + $run1.Length + $run2.Length
$run4End = $run3Start + $run3old.Length + $run4old.Length

$delRange = $d.Range($run3Start, $run4End)
$delRange.Delete()

$insPoint = $d.Range($run3Start, $run3Start)
$insPoint.InsertAfter($run3old + $run4old)

# ---- Step 2: insert the new paragraphs (4, PO, SM, Time) after paragraph 3 ----
$p3b = $d.Paragraphs.Item(3).Range
$insertionPoint = $d.Range($p3b.End, $p3b.End)
$xmlFragment = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t>4)Utilizando o SCRUM, a equipe de projeto seria dívida em:</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>PO(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>Product</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Owner</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>): O Dono da empresa, O PO é o responsável por inspecionar o projeto, acompanhando as entregas participando das reuniões diárias, mostrar as necessidades e garantir o retorno do investimento.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>SM(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ScrumMaster</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>): Eu mesmo, garantindo o uso do Scrum, protege o time das interferências externas além de ser responsável também por resolver eventuais problemas que podem surgir durante o projeto.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Time: Equipe de desenvolvedores, o time é auto gerenciado, uma equipe formada por múltiplas competências, a equipe de desenvolvedores é o time que de fato vai fazer o projeto acontecer, com o </w:t></w:r><w:r><w:t>auxílio</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> do </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ScrumMaster</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> e do PO durante as reuniões diárias. </w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($xmlFragment)

Write-Host "done"
